$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Update the active selection to E8
$ws.Range("E8").Select()
